$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.199.18"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.855.95"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'241.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'0.6988"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.07783"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.3073"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "'23.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").Value = "1.867.84"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "'5.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "'92.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").Value = "'0.6870"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "'6.541"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").Value = "29.196.08"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").Value = "'248.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "2.105.22"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "'12.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'7.542"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").Value = "'0.9998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'0.1504"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "'8.858"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("D28").Value = "'18.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "'1.557"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.02%  "
$ws.Range("D30").Value = "'4.254"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "'4.208"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").Value = "'1.194"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "'0.05234"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "'0.7599"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "'1.844"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").Value = "'2.707"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'0.01861"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").Value = "1.222.03"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "'2.727"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "'0.8990"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").Value = "'109.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").Value = "'0.9994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "2.004.90"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "'5.494"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.51%  "
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").Value = "'65.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.27%  "
$ws.Range("D48").Value = "'0.5183"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").Value = "'9.530"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'1.751"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "'7.052"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.75%  "
